# Apply the updated "cryptos" price/volume snapshot.
# Note: several Price (column D) values look like plain numbers (e.g. "514.99",
# "1.00"), but in this sheet they are plain text. A leading apostrophe is used
# to force Excel to keep them as text instead of silently converting them to
# numeric values (which would drop formatting such as trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.647.15"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").Value = "2.445.27"
$ws.Range("E3").Value = "  -2.16%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'514.99"
$ws.Range("E5").Value = "  -4.17%  "
$ws.Range("D6").Value = "'130.79"
$ws.Range("E6").Value = "  -4.08%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.553"
$ws.Range("E8").Value = "  -2.35%  "
$ws.Range("D9").Value = "2.445.57"
$ws.Range("E9").Value = "  -3.04%  "
$ws.Range("D10").Value = "'0.0975"
$ws.Range("E10").Value = "  -3.65%  "
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("D13").Value = "'0.337"
$ws.Range("E13").Value = "  -3.09%  "
$ws.Range("D14").Value = "2.879.35"
$ws.Range("E14").Value = "  -2.95%  "
$ws.Range("D15").Value = "57.589.21"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D16").Value = "'21.97"
$ws.Range("E16").Value = "  -4.46%  "
$ws.Range("E17").Value = "  -3.39%  "
$ws.Range("D18").Value = "2.484.43"
$ws.Range("E18").Value = "  -1.63%  "
$ws.Range("D19").Value = "'10.57"
$ws.Range("E19").Value = "  -4.73%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'4.13"
$ws.Range("E20").Value = "  -3.29%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'316.16"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("E23").Value = "  -4.32%  "
$ws.Range("D24").Value = "'64.04"
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("D25").Value = "'0.403"
$ws.Range("E25").Value = "  -3.84%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  -3.20%  "
$ws.Range("D28").Value = "'7.26"
$ws.Range("E28").Value = "  -3.53%  "
$ws.Range("E29").Value = "  -5.16%  "
$ws.Range("D30").Value = "'165.12"
$ws.Range("E30").Value = "  -3.61%  "
$ws.Range("D31").Value = "'1.67"
$ws.Range("E31").Value = "  -4.83%  "
$ws.Range("D32").Value = "'6.16"
$ws.Range("E32").Value = "  -6.81%  "
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("D36").Value = "'17.92"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("E37").Value = "  -7.87%  "
$ws.Range("D38").Value = "'3.92"
$ws.Range("E38").Value = "  -3.82%  "
$ws.Range("D39").Value = "'1.45"
$ws.Range("E39").Value = "  -5.46%  "
$ws.Range("D40").Value = "'0.779"
$ws.Range("E40").Value = "  -3.92%  "
$ws.Range("D41").Value = "'3.39"
$ws.Range("E41").Value = "  -5.53%  "
$ws.Range("D42").Value = "'269.60"
$ws.Range("E42").Value = "  -5.37%  "
$ws.Range("D43").Value = "'4.90"
$ws.Range("E43").Value = "  -5.38%  "
$ws.Range("D44").Value = "'0.584"
$ws.Range("E44").Value = "  -3.67%  "
$ws.Range("D45").Value = "'122.88"
$ws.Range("E45").Value = "  -5.73%  "
$ws.Range("E46").Value = "  -2.16%  "
$ws.Range("D47").Value = "'0.0482"
$ws.Range("E47").Value = "  -4.53%  "
$ws.Range("E48").Value = "  -5.51%  "
$ws.Range("D49").Value = "'16.52"
$ws.Range("E49").Value = "  -4.77%  "
$ws.Range("D50").Value = "1.715.14"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("D51").Value = "'0.966"
$ws.Range("E51").Value = "  -2.45%  "
